# Update NATMI LR-pair output (Wnt5b-Fzd2) with newly computed TPM-based
# specificity statistics, and add the new MuSCs-sourced signalling rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Row 2 (FAPs -> Wnt5b -> Fzd2 -> ECs): refreshed specificity values
# ---------------------------------------------------------------------
$ws.Range("I2").Value2  = 0.4348530100317075
$ws.Range("J2").Value2  = 0.4348530100317076
$ws.Range("K2").Value2  = 3
$ws.Range("L2").Value2  = 1
$ws.Range("M2").Value2  = 0.06449866666666666
$ws.Range("N2").Value2  = 0.193496
$ws.Range("O2").Value2  = 0.004525829983623641
$ws.Range("P2").Value2  = 0.004525829983623642
$ws.Range("Q2").Value2  = 0.02609766550222222
$ws.Range("R2").Value2  = 0.23487898952
$ws.Range("S2").Value2  = 0.001968070791270494
$ws.Range("T2").Value2  = 0.001968070791270495

# ---------------------------------------------------------------------
# Row 3 (FAPs -> Wnt5b -> Fzd2 -> FAPs): refreshed specificity values
# ---------------------------------------------------------------------
$ws.Range("I3").Value2  = 0.4348530100317075
$ws.Range("J3").Value2  = 0.4348530100317076
$ws.Range("O3").Value2  = 0.745188142173877
$ws.Range("P3").Value2  = 0.7451881421738772
$ws.Range("S3").Value2  = 0.3240473066642464
$ws.Range("T3").Value2  = 0.3240473066642465

# ---------------------------------------------------------------------
# Row 4 (FAPs -> Wnt5b -> Fzd2 -> MuSCs): refreshed specificity values
# ---------------------------------------------------------------------
$ws.Range("I4").Value2  = 0.4348530100317075
$ws.Range("J4").Value2  = 0.4348530100317076
$ws.Range("M4").Value2  = 3.566885000000001
$ws.Range("N4").Value2  = 10.700655
$ws.Range("O4").Value2  = 0.2502860278424993
$ws.Range("P4").Value2  = 0.2502860278424993
$ws.Range("Q4").Value2  = 1.443244898316667
$ws.Range("R4").Value2  = 12.98920408485
$ws.Range("S4").Value2  = 0.1088376325761906
$ws.Range("T4").Value2  = 0.1088376325761906

# ---------------------------------------------------------------------
# Row 5 (new): MuSCs -> Wnt5b -> Fzd2 -> ECs
# ---------------------------------------------------------------------
$ws.Range("A5").Value2 = "MuSCs"
$ws.Range("B5").Value2 = "Wnt5b"
$ws.Range("C5").Value2 = "Fzd2"
$ws.Range("D5").Value2 = "ECs"
$ws.Range("E5").Value2 = 2
$ws.Range("F5").Value2 = 0.6666666666666666
$ws.Range("G5").Value2 = 0.5258596666666667
$ws.Range("H5").Value2 = 1.577579
$ws.Range("I5").Value2 = 0.5651469899682925
$ws.Range("J5").Value2 = 0.5651469899682925
$ws.Range("K5").Value2 = 3
$ws.Range("L5").Value2 = 1
$ws.Range("M5").Value2 = 0.06449866666666666
$ws.Range("N5").Value2 = 0.193496
$ws.Range("O5").Value2 = 0.004525829983623641
$ws.Range("P5").Value2 = 0.004525829983623642
$ws.Range("Q5").Value2 = 0.03391724735377778
$ws.Range("R5").Value2 = 0.305255226184
$ws.Range("S5").Value2 = 0.002557759192353147
$ws.Range("T5").Value2 = 0.002557759192353148

# ---------------------------------------------------------------------
# Row 6 (new): MuSCs -> Wnt5b -> Fzd2 -> FAPs
# ---------------------------------------------------------------------
$ws.Range("A6").Value2 = "MuSCs"
$ws.Range("B6").Value2 = "Wnt5b"
$ws.Range("C6").Value2 = "Fzd2"
$ws.Range("D6").Value2 = "FAPs"
$ws.Range("E6").Value2 = 2
$ws.Range("F6").Value2 = 0.6666666666666666
$ws.Range("G6").Value2 = 0.5258596666666667
$ws.Range("H6").Value2 = 1.577579
$ws.Range("I6").Value2 = 0.5651469899682925
$ws.Range("J6").Value2 = 0.5651469899682925
$ws.Range("K6").Value2 = 3
$ws.Range("L6").Value2 = 1
$ws.Range("M6").Value2 = 10.61985133333333
$ws.Range("N6").Value2 = 31.859554
$ws.Range("O6").Value2 = 0.745188142173877
$ws.Range("P6").Value2 = 0.7451881421738772
$ws.Range("Q6").Value2 = 5.584551482196223
$ws.Range("R6").Value2 = 50.260963339766
$ws.Range("S6").Value2 = 0.4211408355096306
$ws.Range("T6").Value2 = 0.4211408355096307

# ---------------------------------------------------------------------
# Row 7 (new): MuSCs -> Wnt5b -> Fzd2 -> MuSCs
# ---------------------------------------------------------------------
$ws.Range("A7").Value2 = "MuSCs"
$ws.Range("B7").Value2 = "Wnt5b"
$ws.Range("C7").Value2 = "Fzd2"
$ws.Range("D7").Value2 = "MuSCs"
$ws.Range("E7").Value2 = 2
$ws.Range("F7").Value2 = 0.6666666666666666
$ws.Range("G7").Value2 = 0.5258596666666667
$ws.Range("H7").Value2 = 1.577579
$ws.Range("I7").Value2 = 0.5651469899682925
$ws.Range("J7").Value2 = 0.5651469899682925
$ws.Range("K7").Value2 = 3
$ws.Range("L7").Value2 = 1
$ws.Range("M7").Value2 = 3.566885000000001
$ws.Range("N7").Value2 = 10.700655
$ws.Range("O7").Value2 = 0.2502860278424993
$ws.Range("P7").Value2 = 0.2502860278424993
$ws.Range("Q7").Value2 = 1.875680957138334
$ws.Range("R7").Value2 = 16.881128614245
$ws.Range("S7").Value2 = 0.1414483952663087
$ws.Range("T7").Value2 = 0.1414483952663087
